$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''26.924.79'
$ws.Range('E2').Value = '  -1.66%  '
$ws.Range('D3').Value = '''1.813.38'
$ws.Range('E3').Value = '  -0.65%  '
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').Value = '''310.54'
$ws.Range('E5').Value = '  -0.99%  '
$ws.Range('D6').Value = '''0.9996'
$ws.Range('E6').Value = '  -0.16%  '
$ws.Range('D7').Value = '''0.4641'
$ws.Range('E7').Value = '  +3.82%  '
$ws.Range('D8').Value = '''0.3750'
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('D9').Value = '''0.07430'
$ws.Range('E9').Value = '  -1.05%  '
$ws.Range('D10').Value = '''0.8850'
$ws.Range('E10').Value = '  -0.49%  '
$ws.Range('D11').Value = '''20.49'
$ws.Range('E11').Value = '  -2.63%  '
$ws.Range('D12').Value = '''1.788.14'
$ws.Range('E12').Value = '  -2.09%  '
$ws.Range('D13').Value = '''5.359'
$ws.Range('E13').Value = '  -1.02%  '
$ws.Range('D14').Value = '''6.537'
$ws.Range('E14').Value = '  -3.33%  '
$ws.Range('D15').Value = '''0.07048'
$ws.Range('E15').Value = '  -0.77%  '
$ws.Range('D16').Value = '''91.52'
$ws.Range('E16').Value = '  -2.68%  '
$ws.Range('D17').Value = '''1.001'
$ws.Range('D18').Value = '''0.000008789'
$ws.Range('E18').Value = '  -0.31%  '
$ws.Range('D19').Value = '''0.9991'
$ws.Range('E19').Value = '  -0.19%  '
$ws.Range('D20').Value = '''14.78'
$ws.Range('E20').Value = '  -2.98%  '
$ws.Range('D21').Value = '''26.909.34'
$ws.Range('E21').Value = '  -1.70%  '
$ws.Range('D22').Value = '''5.318'
$ws.Range('E22').Value = '  +1.01%  '
$ws.Range('D23').Value = '''10.80'
$ws.Range('E23').Value = '  -1.12%  '
$ws.Range('D24').Value = '''1.985.23'
$ws.Range('E24').Value = '  -3.39%  '
$ws.Range('D25').Value = '''1.916'
$ws.Range('E25').Value = '  -2.70%  '
$ws.Range('D26').Value = '''151.48'
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('D27').Value = '''18.48'
$ws.Range('E27').Value = '  -0.52%  '
$ws.Range('D28').Value = '''2.168'
$ws.Range('E28').Value = '  -8.90%  '
$ws.Range('D29').Value = '''5.303'
$ws.Range('E29').Value = '  -0.97%  '
$ws.Range('D30').Value = '''115.22'
$ws.Range('E30').Value = '  -2.47%  '
$ws.Range('D31').Value = '''0.08906'
$ws.Range('E31').Value = '  +0.82%  '
$ws.Range('D32').Value = '''0.7718'
$ws.Range('E32').Value = '  -1.81%  '
$ws.Range('E33').Value = '  -2.36%  '
$ws.Range('D34').Value = '''4.490'
$ws.Range('E34').Value = '  -0.63%  '
$ws.Range('D35').Value = '''2.896'
$ws.Range('E35').Value = '  -1.48%  '
$ws.Range('D36').Value = '''0.9995'
$ws.Range('E36').Value = '  -0.13%  '
$ws.Range('E37').Value = '  +0.79%  '
$ws.Range('D38').Value = '''2.485'
$ws.Range('E38').Value = '  +8.14%  '
$ws.Range('D39').Value = '''0.01958'
$ws.Range('E39').Value = '  -1.86%  '
$ws.Range('D40').Value = '''0.05239'
$ws.Range('E40').Value = '  -1.73%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').Value = '''7.222'
$ws.Range('E41').Value = '  -2.15%  '
$ws.Range('B42').Value = 'MXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D42').Value = '''2.909'
$ws.Range('E42').Value = '  +1.78%  '
$ws.Range('D43').Value = '''0.5315'
$ws.Range('E43').Value = '  +0.03%  '
$ws.Range('D44').Value = '''0.1662'
$ws.Range('E44').Value = '  -3.85%  '
$ws.Range('D45').Value = '''8.634'
$ws.Range('E45').Value = '  -1.47%  '
$ws.Range('D46').Value = '''0.5085'
$ws.Range('E46').Value = '  -0.53%  '
$ws.Range('D47').Value = '''10.39'
$ws.Range('E47').Value = '  -2.42%  '
$ws.Range('D48').Value = '''104.61'
$ws.Range('E48').Value = '  -1.15%  '
$ws.Range('D49').Value = '''1.677'
$ws.Range('E49').Value = '  -1.38%  '
$ws.Range('D50').Value = '''0.9990'
$ws.Range('E50').Value = '  -0.22%  '
$ws.Range('D51').Value = '''0.06325'
$ws.Range('E51').Value = '  -0.90%  '
